$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '25.887.03'
Set-TextValue $ws.Range('E2') '  +0.25%  '
Set-TextValue $ws.Range('D3') '1.741.61'
Set-TextValue $ws.Range('E3') '  +0.48%  '
Set-TextValue $ws.Range('E4') '  -0.11%  '
Set-TextValue $ws.Range('D5') '238.55'
Set-TextValue $ws.Range('E5') '  +3.76%  '
Set-TextValue $ws.Range('D6') '0.9998'
Set-TextValue $ws.Range('E6') '  -0.06%  '
Set-TextValue $ws.Range('D7') '0.5168'
Set-TextValue $ws.Range('E7') '  -0.81%  '
Set-TextValue $ws.Range('D8') '0.2739'
Set-TextValue $ws.Range('E8') '  -0.55%  '
Set-TextValue $ws.Range('E9') '  +0.16%  '
Set-TextValue $ws.Range('D10') '1.739.71'
Set-TextValue $ws.Range('E10') '  +0.33%  '
Set-TextValue $ws.Range('D11') '0.07163'
Set-TextValue $ws.Range('E11') '  +1.70%  '
Set-TextValue $ws.Range('E12') '  +1.44%  '
Set-TextValue $ws.Range('D13') '14.96'
Set-TextValue $ws.Range('E13') '  -0.24%  '
Set-TextValue $ws.Range('D14') '4.593'
Set-TextValue $ws.Range('E14') '  +1.57%  '
Set-TextValue $ws.Range('D15') '77.35'
Set-TextValue $ws.Range('E15') '  +1.00%  '
Set-TextValue $ws.Range('D16') '0.9999'
Set-TextValue $ws.Range('E16') '  +0.02%  '
Set-TextValue $ws.Range('D17') '0.9995'
Set-TextValue $ws.Range('E17') '  -0.10%  '
Set-TextValue $ws.Range('D18') '25.896.53'
Set-TextValue $ws.Range('E18') '  +0.28%  '
Set-TextValue $ws.Range('D19') '11.72'
Set-TextValue $ws.Range('D20') '0.000006785'
Set-TextValue $ws.Range('E20') '  +2.32%  '
Set-TextValue $ws.Range('D21') '1.961.70'
Set-TextValue $ws.Range('E21') '  +0.30%  '
Set-TextValue $ws.Range('D22') '4.271'
Set-TextValue $ws.Range('E22') '  +1.94%  '
Set-TextValue $ws.Range('D23') '8.675'
Set-TextValue $ws.Range('E23') '  -1.06%  '
Set-TextValue $ws.Range('D24') '5.245'
Set-TextValue $ws.Range('E24') '  +1.79%  '
Set-TextValue $ws.Range('D25') '138.64'
Set-TextValue $ws.Range('E25') '  -0.75%  '
Set-TextValue $ws.Range('E26') '  +0.35%  '
Set-TextValue $ws.Range('D27') '15.13'
Set-TextValue $ws.Range('E27') '  +0.92%  '
Set-TextValue $ws.Range('D28') '1.763'
Set-TextValue $ws.Range('E28') '  -0.71%  '
Set-TextValue $ws.Range('E29') '  +3.68%  '
Set-TextValue $ws.Range('D30') '3.950'
Set-TextValue $ws.Range('E30') '  +6.12%  '
Set-TextValue $ws.Range('D31') '0.08296'
Set-TextValue $ws.Range('E31') '  +0.04%  '
Set-TextValue $ws.Range('D32') '3.646'
Set-TextValue $ws.Range('E32') '  +4.30%  '
Set-TextValue $ws.Range('D33') '0.04587'
Set-TextValue $ws.Range('E33') '  +2.90%  '
Set-TextValue $ws.Range('D34') '2.661'
Set-TextValue $ws.Range('E34') '  +2.18%  '
Set-TextValue $ws.Range('D35') '0.9893'
Set-TextValue $ws.Range('E35') '  +1.76%  '
Set-TextValue $ws.Range('D36') '0.6183'
Set-TextValue $ws.Range('E36') '  +0.34%  '
Set-TextValue $ws.Range('D37') '2.685'
Set-TextValue $ws.Range('E37') '  +0.72%  '
Set-TextValue $ws.Range('D38') '0.01610'
Set-TextValue $ws.Range('E38') '  +2.59%  '
Set-TextValue $ws.Range('D39') '1.932'
Set-TextValue $ws.Range('E39') '  +2.14%  '
Set-TextValue $ws.Range('D40') '0.9995'
Set-TextValue $ws.Range('D41') '97.85'
Set-TextValue $ws.Range('E41') '  -2.11%  '
Set-TextValue $ws.Range('D42') '0.3836'
Set-TextValue $ws.Range('E42') '  +0.38%  '
Set-TextValue $ws.Range('D43') '0.7391'
Set-TextValue $ws.Range('E43') '  +2.58%  '
Set-TextValue $ws.Range('D44') '4.979'
Set-TextValue $ws.Range('E44') '  -0.67%  '
Set-TextValue $ws.Range('D45') '0.1124'
Set-TextValue $ws.Range('E45') '  +0.14%  '
Set-TextValue $ws.Range('D46') '6.209'
Set-TextValue $ws.Range('E46') '  +0.56%  '
Set-TextValue $ws.Range('D47') '0.05261'
Set-TextValue $ws.Range('E47') '  -1.47%  '
Set-TextValue $ws.Range('D48') '54.83'
Set-TextValue $ws.Range('E48') '  +3.10%  '
Set-TextValue $ws.Range('D49') '30.51'
Set-TextValue $ws.Range('E49') '  +1.89%  '
Set-TextValue $ws.Range('D50') '7.579'
Set-TextValue $ws.Range('E50') '  -0.71%  '
Set-TextValue $ws.Range('E51') '  +0.99%  '
